# Teilaufgaben.xlsx -- "gui ueberarbeitet, thread geaendert (laeuft
# stabiler), neue befehle implementiert, interrupts implementiert ..."
#
# Updates the "Was laeuft durch?" test-matrix on Tabelle1 with a few more
# marked/unmarked test results and moves the view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Tabelle1 (tabSelected="1")

# --- G5 now also gets the "x" mark, like G3/G4/G8.
# The cell was an empty, styled placeholder (s="3"); clear its format first
# so the new value lands unstyled, same as the authored edit.
$ws.Range("G5").ClearFormats()
$ws.Range("G5").Value = "x"

# --- New entries in column C (points awarded per test row).
# "------------" must become shared-string index 35 and "???" index 36,
# so write the "------------" cell before any "???" cell.
$ws.Range("C13").Value = 4
$ws.Range("C14").Value = "'------------"
$ws.Range("C9").Value = "???"
$ws.Range("C15").Value = "???"
$ws.Range("C16").Value = 5
$ws.Range("C26").Value = 5

# --- View: scroll the window down and move the active selection to C24.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select()
